$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.639.25"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "2.473.66"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.51"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.58"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.00"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0852"
$ws.Range("E11").Value = "  +8.23%  "
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "2.854.95"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.72"
$ws.Range("E15").Value = "  -3.41%  "
$ws.Range("D16").Value = "2.479.95"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("D18").Value = "41.598.90"
$ws.Range("D19").Value = "0.0₃0950"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.28"
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.28"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.82"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.18"
$ws.Range("E31").Value = "  +1.95%  "
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0767"
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.27"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.86"
$ws.Range("E37").Value = "  +2.55%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.92"
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.115"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("E42").Value = "  +3.71%  "
$ws.Range("D43").Value = "1.986.09"
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.74"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("E46").Value = "  +1.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.28"
$ws.Range("E47").Value = "  +2.67%  "
$ws.Range("D48").Value = "2.712.48"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.15"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.04"
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.03"
$ws.Range("E51").Value = "  -0.41%  "
